$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16: GT vs LSG
$ws.Range("E16").Value = 70
$ws.Range("H16").Value = 50
$ws.Range("K16").Value = 100
$ws.Range("N16").Value = 0
$ws.Range("Q16").Value = 40
$ws.Range("T16").Value = 80
$ws.Range("W16").Value = 60

# Row 17: SRH vs RR (Contest 5)
$ws.Range("E17").Value = 0
$ws.Range("H17").Value = 50
$ws.Range("K17").Value = 60
$ws.Range("N17").Value = 70
$ws.Range("Q17").Value = 100
$ws.Range("T17").Value = 40
$ws.Range("W17").Value = 80
